# Lang TODO.xlsx -- "Naprawione edytowanie opcji zapytan." edit
#
# 1. Add a new backlog row (#14) describing the QueryOptions/[Displayed]
#    fallback-selection bug.
# 2. Row 10 no longer needs its custom (wrapped-text) height.
# 3. Conditional formatting ranges grow to cover the new row.
# 4. Selection lands on A16 (the row right below the new data), matching
#    where Excel leaves the cursor after typing the new row and hitting Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 15: Id=14, Zadanie=<new text>, Moduł=logika ---------------
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Jeżeli dla danego pytania nie ma ani jednego QueryOptions z [Displayed] = True, wtedy ma wybierać spośród pozostałych."
$ws.Cells.Item(15, 4).Value = "logika"
$ws.Rows.Item(15).RowHeight = 25.5

# --- Row 10 goes back to the default row height -------------------------
$ws.Rows.Item(10).AutoFit()

# --- Conditional formatting now spans through the new row (+ a little
#     headroom, matching the authored sqref of A2:F17 / F10:F17) ---------
$ws.Cells.FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A2:F17"))
$ws.Cells.FormatConditions.Item(4).ModifyAppliesToRange($ws.Range("F10:F17"))

# --- Leave the selection where Excel would after entering the new row ---
$ws.Range("A16").Select()
